$wb = $excel.ActiveWorkbook

# Rename worksheets
$wb.Worksheets.Item(1).Name = "DT vs veh"
$wb.Worksheets.Item(2).Name = "DT+Treg vs veh"
$wb.Worksheets.Item(3).Name = "DT+Treg vs DT"

# --- Worksheet 1 (sheet1) updates ---
$ws = $wb.Worksheets.Item(1)
# Row 2
$ws.Cells.Item(2,5).Value = "'6037"
$ws.Cells.Item(2,6).Value = [double]"1.917565574867371e-05"
$ws.Cells.Item(2,7).Value = [double]"0.004508384029256892"
$ws.Cells.Item(2,8).Value = [double]"0.004199403525926525"
$ws.Cells.Item(2,11).Value = [double]"2.345979097354497"
# Row 3
$ws.Cells.Item(3,5).Value = "'6037"
$ws.Cells.Item(3,6).Value = [double]"3.243441747666829e-05"
$ws.Cells.Item(3,7).Value = [double]"0.004508384029256892"
$ws.Cells.Item(3,8).Value = [double]"0.004199403525926525"
$ws.Cells.Item(3,11).Value = [double]"2.345979097354497"
# Row 4
$ws.Cells.Item(4,5).Value = "'6037"
$ws.Cells.Item(4,6).Value = [double]"8.226084278446349e-05"
$ws.Cells.Item(4,7).Value = [double]"0.00762283809802695"
$ws.Cells.Item(4,8).Value = [double]"0.007100409587711586"
$ws.Cells.Item(4,11).Value = [double]"2.117883304141839"
# Row 5
$ws.Cells.Item(5,5).Value = "'6037"
$ws.Cells.Item(5,6).Value = [double]"0.0001799759521517347"
$ws.Cells.Item(5,7).Value = [double]"0.009225776087286945"
$ws.Cells.Item(5,8).Value = [double]"0.008593490789369893"
$ws.Cells.Item(5,11).Value = [double]"2.03499709007414"
# Row 6
$ws.Cells.Item(6,5).Value = "'6037"
$ws.Cells.Item(6,6).Value = [double]"0.0001799759521517347"
$ws.Cells.Item(6,7).Value = [double]"0.009225776087286945"
$ws.Cells.Item(6,8).Value = [double]"0.008593490789369893"
$ws.Cells.Item(6,11).Value = [double]"2.03499709007414"
# Row 7
$ws.Cells.Item(7,5).Value = "'6037"
$ws.Cells.Item(7,6).Value = [double]"0.0002192103794309596"
$ws.Cells.Item(7,7).Value = [double]"0.009225776087286945"
$ws.Cells.Item(7,8).Value = [double]"0.008593490789369893"
$ws.Cells.Item(7,11).Value = [double]"2.03499709007414"
# Row 8
$ws.Cells.Item(8,5).Value = "'6037"
$ws.Cells.Item(8,6).Value = [double]"0.0002323037144280885"
$ws.Cells.Item(8,7).Value = [double]"0.009225776087286945"
$ws.Cells.Item(8,8).Value = [double]"0.008593490789369893"
$ws.Cells.Item(8,11).Value = [double]"2.03499709007414"
# Row 9
$ws.Cells.Item(9,5).Value = "'6037"
$ws.Cells.Item(9,6).Value = [double]"0.0003724083919974529"
$ws.Cells.Item(9,7).Value = [double]"0.01294119162191149"
$ws.Cells.Item(9,8).Value = [double]"0.01205427163570703"
$ws.Cells.Item(9,11).Value = [double]"1.888025732091304"
# Row 10
$ws.Cells.Item(10,5).Value = "'6037"
$ws.Cells.Item(10,6).Value = [double]"0.0007306881056198543"
$ws.Cells.Item(10,7).Value = [double]"0.02257014370692439"
$ws.Cells.Item(10,8).Value = [double]"0.02102330689853616"
$ws.Cells.Item(10,11).Value = [double]"1.646465675706217"
# Row 11
$ws.Cells.Item(11,5).Value = "'6037"
$ws.Cells.Item(11,6).Value = [double]"0.0008495549940649225"
$ws.Cells.Item(11,7).Value = [double]"0.02293854243485529"
$ws.Cells.Item(11,8).Value = [double]"0.02136645755007347"
$ws.Cells.Item(11,11).Value = [double]"1.639434181579059"
# Row 12
$ws.Cells.Item(12,5).Value = "'6037"
$ws.Cells.Item(12,6).Value = [double]"0.0009733938302646022"
$ws.Cells.Item(12,7).Value = [double]"0.02293854243485529"
$ws.Cells.Item(12,8).Value = [double]"0.02136645755007347"
$ws.Cells.Item(12,11).Value = [double]"1.639434181579059"
# Row 13
$ws.Cells.Item(13,5).Value = "'6037"
$ws.Cells.Item(13,6).Value = [double]"0.000990152910857063"
$ws.Cells.Item(13,7).Value = [double]"0.02293854243485529"
$ws.Cells.Item(13,8).Value = [double]"0.02136645755007347"
$ws.Cells.Item(13,11).Value = [double]"1.639434181579059"
# Row 14
$ws.Cells.Item(14,5).Value = "'6037"
$ws.Cells.Item(14,6).Value = [double]"0.001221558016871538"
$ws.Cells.Item(14,7).Value = [double]"0.02612254836079135"
$ws.Cells.Item(14,8).Value = [double]"0.02433224875711728"
$ws.Cells.Item(14,11).Value = [double]"1.58298445813799"
# Row 15
$ws.Cells.Item(15,5).Value = "'6037"
$ws.Cells.Item(15,6).Value = [double]"0.002313137282654205"
$ws.Cells.Item(15,7).Value = [double]"0.04593229746984778"
$ws.Cells.Item(15,8).Value = [double]"0.04278434372428078"
$ws.Cells.Item(15,11).Value = [double]"1.337881831230823"
# Row 16
$ws.Cells.Item(16,5).Value = "'6037"
$ws.Cells.Item(16,6).Value = [double]"0.002983429875789102"
$ws.Cells.Item(16,7).Value = [double]"0.0552929003646247"
$ws.Cells.Item(16,8).Value = [double]"0.05150342101362239"
$ws.Cells.Item(16,11).Value = [double]"1.257330628738869"
# Row 17
$ws.Cells.Item(17,5).Value = "'6037"
$ws.Cells.Item(17,6).Value = [double]"0.004315716288945614"
$ws.Cells.Item(17,7).Value = [double]"0.07498557052043003"
$ws.Cells.Item(17,8).Value = [double]"0.06984646099214611"
$ws.Cells.Item(17,11).Value = [double]"1.125022299891775"
# Row 18
$ws.Cells.Item(18,5).Value = "'6037"
$ws.Cells.Item(18,6).Value = [double]"0.006065389928469733"
$ws.Cells.Item(18,7).Value = [double]"0.09817852579322092"
$ws.Cells.Item(18,8).Value = [double]"0.09144989528637769"
$ws.Cells.Item(18,11).Value = [double]"1.007983493367501"
# Row 19
$ws.Cells.Item(19,5).Value = "'6037"
$ws.Cells.Item(19,6).Value = [double]"0.006356882965028693"
$ws.Cells.Item(19,7).Value = [double]"0.09817852579322092"
$ws.Cells.Item(19,8).Value = [double]"0.09144989528637769"
$ws.Cells.Item(19,11).Value = [double]"1.007983493367501"
# Row 20
$ws.Cells.Item(20,5).Value = "'6037"
$ws.Cells.Item(20,6).Value = [double]"0.006733883476979991"
$ws.Cells.Item(20,7).Value = [double]"0.09852734771581249"
$ws.Cells.Item(20,8).Value = [double]"0.09177481082199877"
$ws.Cells.Item(20,11).Value = [double]"1.00644320794295"
# Row 21
$ws.Cells.Item(21,5).Value = "'6037"
$ws.Cells.Item(21,6).Value = [double]"0.008449346044309991"
$ws.Cells.Item(21,7).Value = [double]"0.1174459100159089"
$ws.Cells.Item(21,8).Value = [double]"0.1093967961526452"
$ws.Cells.Item(21,11).Value = [double]"0.9301621026705392"
# Row 22
$ws.Cells.Item(22,5).Value = "'6037"
$ws.Cells.Item(22,6).Value = [double]"0.009247341709563046"
$ws.Cells.Item(22,7).Value = [double]"0.1217087938063237"
$ws.Cells.Item(22,8).Value = [double]"0.1133675247116836"
$ws.Cells.Item(22,11).Value = [double]"0.9146780416254197"
# Row 23
$ws.Cells.Item(23,5).Value = "'6037"
$ws.Cells.Item(23,6).Value = [double]"0.009631631164529213"
$ws.Cells.Item(23,7).Value = [double]"0.1217087938063237"
$ws.Cells.Item(23,8).Value = [double]"0.1133675247116836"
$ws.Cells.Item(23,11).Value = [double]"0.9146780416254197"
# Row 24
$ws.Cells.Item(24,5).Value = "'6037"
$ws.Cells.Item(24,6).Value = [double]"0.01023698959876024"
$ws.Cells.Item(24,7).Value = [double]"0.1237340481937108"
$ws.Cells.Item(24,8).Value = [double]"0.1152539790066371"
$ws.Cells.Item(24,11).Value = [double]"0.9075107780751513"
# Row 25
$ws.Cells.Item(25,5).Value = "'6037"
$ws.Cells.Item(25,6).Value = [double]"0.011374892342548"
$ws.Cells.Item(25,7).Value = [double]"0.1317591696345144"
$ws.Cells.Item(25,8).Value = [double]"0.1227291015906495"
$ws.Cells.Item(25,11).Value = [double]"0.8802191508167835"
# Row 26
$ws.Cells.Item(26,5).Value = "'6037"
$ws.Cells.Item(26,6).Value = [double]"0.01303854470760689"
$ws.Cells.Item(26,7).Value = [double]"0.1449886171485886"
$ws.Cells.Item(26,8).Value = [double]"0.1350518736030019"
$ws.Cells.Item(26,11).Value = [double]"0.838666092272647"

# --- Worksheet 2 (sheet2) updates ---
$ws = $wb.Worksheets.Item(2)
# Row 2
$ws.Cells.Item(2,5).Value = "'6037"
$ws.Cells.Item(2,6).Value = [double]"0.0001526114206987239"
$ws.Cells.Item(2,7).Value = [double]"0.04746215183730313"
$ws.Cells.Item(2,8).Value = [double]"0.04722921861623666"
$ws.Cells.Item(2,11).Value = [double]"1.323652575644986"

# --- Worksheet 3 (sheet3) updates ---
$ws = $wb.Worksheets.Item(3)
# Row 2
$ws.Cells.Item(2,5).Value = "'6037"
$ws.Cells.Item(2,6).Value = [double]"5.096696880434185e-06"
$ws.Cells.Item(2,7).Value = [double]"0.0006319904131738389"
$ws.Cells.Item(2,8).Value = [double]"0.0005794139611440968"
$ws.Cells.Item(2,11).Value = [double]"3.199289509593057"
# Row 3
$ws.Cells.Item(3,5).Value = "'6037"
$ws.Cells.Item(3,6).Value = [double]"0.001055759435675504"
$ws.Cells.Item(3,7).Value = [double]"0.0311389916748292"
$ws.Cells.Item(3,8).Value = [double]"0.02854848133176191"
$ws.Cells.Item(3,11).Value = [double]"1.506695454617855"
# Row 4
$ws.Cells.Item(4,5).Value = "'6037"
$ws.Cells.Item(4,6).Value = [double]"0.001055759435675504"
$ws.Cells.Item(4,7).Value = [double]"0.0311389916748292"
$ws.Cells.Item(4,8).Value = [double]"0.02854848133176191"
$ws.Cells.Item(4,11).Value = [double]"1.506695454617855"
# Row 5
$ws.Cells.Item(5,5).Value = "'6037"
$ws.Cells.Item(5,6).Value = [double]"0.001198979649066336"
$ws.Cells.Item(5,7).Value = [double]"0.0311389916748292"
$ws.Cells.Item(5,8).Value = [double]"0.02854848133176191"
$ws.Cells.Item(5,11).Value = [double]"1.506695454617855"
# Row 6
$ws.Cells.Item(6,5).Value = "'6037"
$ws.Cells.Item(6,6).Value = [double]"0.001353894711080967"
$ws.Cells.Item(6,7).Value = [double]"0.0311389916748292"
$ws.Cells.Item(6,8).Value = [double]"0.02854848133176191"
$ws.Cells.Item(6,11).Value = [double]"1.506695454617855"
# Row 7
$ws.Cells.Item(7,5).Value = "'6037"
$ws.Cells.Item(7,6).Value = [double]"0.001520852224436363"
$ws.Cells.Item(7,7).Value = [double]"0.0311389916748292"
$ws.Cells.Item(7,8).Value = [double]"0.02854848133176191"
$ws.Cells.Item(7,11).Value = [double]"1.506695454617855"
# Row 8
$ws.Cells.Item(8,5).Value = "'6037"
$ws.Cells.Item(8,6).Value = [double]"0.001757846304224229"
$ws.Cells.Item(8,7).Value = [double]"0.0311389916748292"
$ws.Cells.Item(8,8).Value = [double]"0.02854848133176191"
$ws.Cells.Item(8,11).Value = [double]"1.506695454617855"
# Row 9
$ws.Cells.Item(9,5).Value = "'6037"
$ws.Cells.Item(9,6).Value = [double]"0.004192977833011918"
$ws.Cells.Item(9,7).Value = [double]"0.0635183901023383"
$ws.Cells.Item(9,8).Value = [double]"0.05823417768295871"
$ws.Cells.Item(9,11).Value = [double]"1.19710051779965"
# Row 10
$ws.Cells.Item(10,5).Value = "'6037"
$ws.Cells.Item(10,6).Value = [double]"0.004610205733234231"
$ws.Cells.Item(10,7).Value = [double]"0.0635183901023383"
$ws.Cells.Item(10,8).Value = [double]"0.05823417768295871"
$ws.Cells.Item(10,11).Value = [double]"1.19710051779965"
# Row 11
$ws.Cells.Item(11,5).Value = "'6037"
$ws.Cells.Item(11,6).Value = [double]"0.0073559655591904"
$ws.Cells.Item(11,7).Value = [double]"0.09121397293396095"
$ws.Cells.Item(11,8).Value = [double]"0.08362571372553296"
$ws.Cells.Item(11,11).Value = [double]"1.039938627644289"
# Row 12
$ws.Cells.Item(12,5).Value = "'6037"
$ws.Cells.Item(12,6).Value = [double]"0.01278857342509136"
$ws.Cells.Item(12,7).Value = [double]"0.1441621004283026"
$ws.Cells.Item(12,8).Value = [double]"0.1321689885081212"
$ws.Cells.Item(12,11).Value = [double]"0.8411488986889852"
